# Grade sheet updated with comments and scores for all four assignments

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Assignment 1
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Assignment 1")
$ws1.Range("B14").Value = 8

# ---------------------------------------------------------------------------
# Assignment 2
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Assignment 2")
$ws2.Columns.Item(3).ColumnWidth = 68.66666666666667

$ws2.Range("B29").Value = 9.5
$ws2.Range("C29").Value = "Correct but incomplete. Did not show java class file/jar file creation process"

$ws2.Range("B33").Value = 9.5
$ws2.Range("C33").Value = "Submission has no audio"

# ---------------------------------------------------------------------------
# Assignment 3
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("Assignment 3")
$ws3.Columns.Item(3).ColumnWidth = 51.33333333333333

$ws3.Range("B4").Value = 9
$ws3.Range("C4").Value = "No Video submitted"

$ws3.Range("B5").Value = 7
$ws3.Range("C5").Value = "Questions 3.2 and 3.4 are incomplete and are partially incorrect. Did not submit the code."

$ws3.Range("C9").Value = "Answers for questions 3.2 and 3.4 are partially incorrect."

$ws3.Range("C13").Value = "Complete and correct. Did not submit the code."

$ws3.Range("B15").Value = 9
$ws3.Range("C15").Value = "Did not submit the code"

$ws3.Range("C20").Value = "Answer for question 3.3 is partially incorrect (Found least wins when the question is to find the most wins)"

$ws3.Range("C21").Value = "Questions 3.2 and 3.4 are incomplete and are partially incorrect."

$ws3.Range("B29").Value = 9
$ws3.Range("C29").Value = "Did not submit the code"

$ws3.Range("B33").Value = 7.5
$ws3.Range("C33").Value = "Question 3.4's execution was not shown in the video. Question 3.2 was partiaally incorrect. Did not submit code. Video submission had no audio."

$ws3.Range("C34").Value = "Answers for questions 3.2 and 3.4 are partially incorrect."

# ---------------------------------------------------------------------------
# Assignment 4
# ---------------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("Assignment 4")

$ws4.Range("C22").Value = "Complete and correct (Explained the results from a previous session. Did not show the execution of any queries)"

$ws4.Range("B33").Value = 9.5
$ws4.Range("C33").Value = "Submission has no audio."
